$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "all": add today's (2020-05-01 serial 43961) row as new row 33,
# pushing the footnote row down to row 34.
# ---------------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("all")
$wsAll.Activate()

$wsAll.Rows("33:33").Insert()
$wsAll.Range("A33").Value = 43961
$wsAll.Range("B33").Value = 277
$wsAll.Range("C33").Value = 272
$wsAll.Range("D33").Value = 77
$wsAll.Range("E33").Value = 67
$wsAll.Range("F33").Value = 10
$wsAll.Range("G33").Value = 8
$wsAll.Range("H33").Value = 187
$wsAll.Range("H34").Select()

# ---------------------------------------------------------------------------
# Sheet "kobe": yesterday's row (87) gets revised totals, then add today's
# row as new row 88, pushing the footnote row down to row 89.
# ---------------------------------------------------------------------------
$wsKobe = $wb.Worksheets.Item("kobe")
$wsKobe.Activate()

$wsKobe.Range("D87").Value = 1
$wsKobe.Range("E87").Value = 277

$wsKobe.Rows("88:88").Insert()
$wsKobe.Range("A88").Value = 43961
$wsKobe.Range("B88").Value = 0
$wsKobe.Range("C88").Value = 2562
$wsKobe.Range("D88").Value = 0
$wsKobe.Range("E88").Value = 277
$wsKobe.Range("F88").Value = 72
$wsKobe.Range("G88").Value = 63
$wsKobe.Range("H88").Value = 9
$wsKobe.Range("I88").Value = 8
$wsKobe.Range("J88").Value = 178
$wsKobe.Range("K88").Select()

# ---------------------------------------------------------------------------
# Sheet "other": add today's row as new row 63, pushing the footnote row
# down to row 64.
# ---------------------------------------------------------------------------
$wsOther = $wb.Worksheets.Item("other")
$wsOther.Activate()

$wsOther.Rows("63:63").Insert()
$wsOther.Range("A63").Value = 43961
$wsOther.Range("B63").Value = 0
$wsOther.Range("C63").Value = 14
$wsOther.Range("D63").Value = 5
$wsOther.Range("E63").Value = 4
$wsOther.Range("F63").Value = 1
$wsOther.Range("G63").Value = 0
$wsOther.Range("H63").Value = 9
$wsOther.Range("A63").Select()
